$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1435.028035333333
$ws.Range("H2").Value = 4305.084106
$ws.Range("I2").Value = 0.3003423969824614
$ws.Range("J2").Value = 0.3003423969824615
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09976533333333333
$ws.Range("N2").Value = 0.299296
$ws.Range("O2").Value = 0.5779577523264414
$ws.Range("P2").Value = 0.5779577523264414
$ws.Range("Q2").Value = 143.1660502877085
$ws.Range("R2").Value = 1288.494452589376
$ws.Range("S2").Value = 0.1735852166883192
$ws.Range("T2").Value = 0.1735852166883192

# Row 3
$ws.Range("G3").Value = 1435.028035333333
$ws.Range("H3").Value = 4305.084106
$ws.Range("I3").Value = 0.3003423969824614
$ws.Range("J3").Value = 0.3003423969824615
$ws.Range("O3").Value = 0.3247864733292009
$ws.Range("P3").Value = 0.3247864733292009
$ws.Range("Q3").Value = 80.45293343024957
$ws.Range("R3").Value = 724.076400872246
$ws.Range("S3").Value = 0.09754714790717248
$ws.Range("T3").Value = 0.09754714790717249

# Row 4
$ws.Range("G4").Value = 1435.028035333333
$ws.Range("H4").Value = 4305.084106
$ws.Range("I4").Value = 0.3003423969824614
$ws.Range("J4").Value = 0.3003423969824615
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016788
$ws.Range("N4").Value = 0.050364
$ws.Range("O4").Value = 0.09725577434435775
$ws.Range("P4").Value = 0.09725577434435775
$ws.Range("Q4").Value = 24.091250657176
$ws.Range("R4").Value = 216.821255914584
$ws.Range("S4").Value = 0.02921003238696978
$ws.Range("T4").Value = 0.02921003238696979

# Row 5
$ws.Range("I5").Value = 0.1429577688896155
$ws.Range("J5").Value = 0.1429577688896155
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09976533333333333
$ws.Range("N5").Value = 0.299296
$ws.Range("O5").Value = 0.5779577523264414
$ws.Range("P5").Value = 0.5779577523264414
$ws.Range("Q5").Value = 68.1445554656889
$ws.Range("R5").Value = 613.3009991912
$ws.Range("S5").Value = 0.08262355078504503
$ws.Range("T5").Value = 0.08262355078504505

# Row 6
$ws.Range("I6").Value = 0.1429577688896155
$ws.Range("J6").Value = 0.1429577688896155
$ws.Range("O6").Value = 0.3247864733292009
$ws.Range("P6").Value = 0.3247864733292009
$ws.Range("S6").Value = 0.04643074959266916
$ws.Range("T6").Value = 0.04643074959266917

# Row 7
$ws.Range("I7").Value = 0.1429577688896155
$ws.Range("J7").Value = 0.1429577688896155
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.016788
$ws.Range("N7").Value = 0.050364
$ws.Range("O7").Value = 0.09725577434435775
$ws.Range("P7").Value = 0.09725577434435775
$ws.Range("Q7").Value = 11.4670172387
$ws.Range("R7").Value = 103.2031551483
$ws.Range("S7").Value = 0.01390346851190129
$ws.Range("T7").Value = 0.01390346851190129

# Row 8
$ws.Range("G8").Value = 247.7155763333333
$ws.Range("H8").Value = 743.1467289999999
$ws.Range("I8").Value = 0.05184532157837839
$ws.Range("J8").Value = 0.05184532157837839
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.09976533333333333
$ws.Range("N8").Value = 0.299296
$ws.Range("O8").Value = 0.5779577523264414
$ws.Range("P8").Value = 0.5779577523264414
$ws.Range("Q8").Value = 24.71342704475378
$ws.Range("R8").Value = 222.420843402784
$ws.Range("S8").Value = 0.02996440552808112
$ws.Range("T8").Value = 0.02996440552808113

# Row 9
$ws.Range("G9").Value = 247.7155763333333
$ws.Range("H9").Value = 743.1467289999999
$ws.Range("I9").Value = 0.05184532157837839
$ws.Range("J9").Value = 0.05184532157837839
$ws.Range("O9").Value = 0.3247864733292009
$ws.Range("P9").Value = 0.3247864733292009
$ws.Range("Q9").Value = 13.88784349969322
$ws.Range("R9").Value = 124.990591497239
$ws.Range("S9").Value = 0.01683865915405984
$ws.Range("T9").Value = 0.01683865915405984

# Row 10
$ws.Range("G10").Value = 247.7155763333333
$ws.Range("H10").Value = 743.1467289999999
$ws.Range("I10").Value = 0.05184532157837839
$ws.Range("J10").Value = 0.05184532157837839
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.016788
$ws.Range("N10").Value = 0.050364
$ws.Range("O10").Value = 0.09725577434435775
$ws.Range("P10").Value = 0.09725577434435775
$ws.Range("Q10").Value = 4.158649095484
$ws.Range("R10").Value = 37.42784185935599
$ws.Range("S10").Value = 0.00504225689623743
$ws.Range("T10").Value = 0.005042256896237431

# Row 11
$ws.Range("G11").Value = 2412.181518666667
$ws.Range("H11").Value = 7236.544556000001
$ws.Range("I11").Value = 0.5048545125495446
$ws.Range("J11").Value = 0.5048545125495447
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.09976533333333333
$ws.Range("N11").Value = 0.299296
$ws.Range("O11").Value = 0.5779577523264414
$ws.Range("P11").Value = 0.5779577523264414
$ws.Range("Q11").Value = 240.6520932702862
$ws.Range("R11").Value = 2165.868839432576
$ws.Range("S11").Value = 0.291784579324996
$ws.Range("T11").Value = 0.2917845793249961

# Row 12
$ws.Range("G12").Value = 2412.181518666667
$ws.Range("H12").Value = 7236.544556000001
$ws.Range("I12").Value = 0.5048545125495446
$ws.Range("J12").Value = 0.5048545125495447
$ws.Range("O12").Value = 0.3247864733292009
$ws.Range("P12").Value = 0.3247864733292009
$ws.Range("Q12").Value = 135.2357406020218
$ws.Range("R12").Value = 1217.121665418196
$ws.Range("S12").Value = 0.1639699166752994
$ws.Range("T12").Value = 0.1639699166752994

# Row 13
$ws.Range("G13").Value = 2412.181518666667
$ws.Range("H13").Value = 7236.544556000001
$ws.Range("I13").Value = 0.5048545125495446
$ws.Range("J13").Value = 0.5048545125495447
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.016788
$ws.Range("N13").Value = 0.050364
$ws.Range("O13").Value = 0.09725577434435775
$ws.Range("P13").Value = 0.09725577434435775
$ws.Range("Q13").Value = 40.49570333537601
$ws.Range("R13").Value = 364.461330018384
$ws.Range("S13").Value = 0.04910001654924924
$ws.Range("T13").Value = 0.04910001654924925

